$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data fixes / fresh data ---

# New id column (A2)
$ws.Range("A2").Value = 103

# name (B2) gets the Arabic place name plus a tab-separated French label
$ws.Range("B2").Value = "ساحة سلفادور ألندي `tPlace Salvador Allende"

# locale_4 (I2) duplicate of the name is removed entirely (bug fix)
$ws.Range("I2").ClearContents()

# New google_maps_link (X2) with hyperlink, formatted like the other link cells
$ws.Hyperlinks.Add($ws.Range("X2"), "https://goo.gl/maps/9QP1kuoswBG4RZ1s9")
$ws.Range("X2").Style = "Hyperlink"

# New oldest_known_source (Q2) value
$ws.Range("Q2").Value = "abacq date posted"

# Move the visible selection cursor
$ws.Range("Q4").Select()
